$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "2019" column (AE) mirrors the existing year header formatting (text,
# not auto-converted to a number) used by B1:AD1.
$ws.Range("AE1").NumberFormat = "@"
$ws.Range("AE1").Value = "2019"

# Data rows 2-11 for the new 2019 column.
$ws.Range("AE2").Value = 7
$ws.Range("AE3").Value = 2
$ws.Range("AE4").Value = 25
$ws.Range("AE5").Value = 8
$ws.Range("AE6").Value = 2
$ws.Range("AE7").Value = 5
$ws.Range("AE8").Value = 34
$ws.Range("AE9").Value = 7
$ws.Range("AE10").Value = 49
$ws.Range("AE11").Value = 14
